$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values for rows 2-11
$ws.Range("C2").Value = 1.100809410922692
$ws.Range("D2").Value = 0.2828787030002022

$ws.Range("C3").Value = -1.16569513054456
$ws.Range("D3").Value = 0.256220747945666

$ws.Range("C4").Value = -0.9096851982215011
$ws.Range("D4").Value = 0.3728476226108648

$ws.Range("C5").Value = 0.1374289150963754
$ws.Range("D5").Value = 0.891941794213305

$ws.Range("C6").Value = -1.493541334338836
$ws.Range("D6").Value = 0.1495000251924989

$ws.Range("C7").Value = -1.165727686716384
$ws.Range("D7").Value = 0.2562078564913526

$ws.Range("C8").Value = -0.723617842647205
$ws.Range("D8").Value = 0.4769247564186445

$ws.Range("C9").Value = 0.2995676489021638
$ws.Range("D9").Value = 0.7673189114072136

$ws.Range("C10").Value = 1.19693178026341
$ws.Range("D10").Value = 0.244072370699715

$ws.Range("C11").Value = 0.7630321345303647
$ws.Range("D11").Value = 0.4535490904393793
